$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.949.53"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "1.637.44"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  +0.96%  "

$ws.Range("D5").Value = "'214.68"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("E7").Value = "  +0.88%  "

$ws.Range("E8").Value = "  -0.79%  "

$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("D10").Value = "'19.63"
$ws.Range("E10").Value = "  -0.59%  "

$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("D12").Value = "1.864.31"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.25"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.626.89"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("E15").Value = "  -1.32%  "

$ws.Range("D16").Value = "0.0₃0757"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("D17").Value = "'62.49"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").Value = "25.951.33"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D20").Value = "'193.61"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").Value = "'4.37"
$ws.Range("E21").Value = "  -1.42%  "

$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("E23").Value = "  -1.43%  "

$ws.Range("D24").Value = "'144.22"
$ws.Range("E24").Value = "  +1.86%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("E26").Value = "  +0.96%  "

$ws.Range("E27").Value = "  +2.45%  "

$ws.Range("E28").Value = "  -0.38%  "

$ws.Range("D29").Value = "'15.47"

$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").Value = "'0.0502"
$ws.Range("E31").Value = "  +1.60%  "

$ws.Range("E32").Value = "  -1.11%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("E34").Value = "  -2.78%  "

$ws.Range("E35").Value = "  +1.85%  "

$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("D37").Value = "1.138.90"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").Value = "'0.545"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("E40").Value = "  +0.80%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'99.35"
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.798"
$ws.Range("E43").Value = "  +0.43%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.41"
$ws.Range("E44").Value = "  -2.58%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.773.63"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0114"
$ws.Range("E46").Value = "  +8.65%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'56.44"
$ws.Range("E47").Value = "  +1.53%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0530"
$ws.Range("E48").Value = "  +2.86%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.46"
$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.66"
$ws.Range("E50").Value = "  +0.99%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.415"
$ws.Range("E51").Value = "  +0.13%  "
